# Delta - Travel Update Center Page - 20 Test Cases (COMPLETE)
# Adds 6 new worksheets (SACRU - URLs, OUATT - Header Names, OUATT - DeltaSkyClub URL,
# OUATT - FlyDeltaApp URL, OUATT - FlyingPartnerAir URL, FAQ Header Names) at the end of
# the workbook, populates them, and moves the active/selected tab to the new last sheet.

$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd($name) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
    $ws.Name = $name
    return $ws
}

# ---------------------------------------------------------------------------
# 1) SACRU - URLs
# ---------------------------------------------------------------------------
$sacru = Add-SheetAtEnd "SACRU - URLs"
$sacru.Range("A1").Value = "https://www.delta.com/us/en/travel-update-center/extending-skymiles-benefits?src=benefits3"
$sacru.Range("A2").Value = "http://amex.co/extrasupport?src=amex3"
$sacru.Range("A3").Value = "https://www.delta.com/us/en/travel-update-center/ways-we-are-keeping-you-safe/onboard-services?src=obsspace1#socialdist"
$sacru.Range("A4").Value = "https://www.delta.com/us/en/skymiles/program-resources/news-and-updates?src=sm1a#recent"
$sacru.Range("A5").Value = "https://www.delta.com/us/en/travel-update-center/extending-skymiles-benefits?src=benefits1"
$sacru.Range("A6").Value = "https://www.delta.com/us/en/travel-update-center/ways-we-are-keeping-you-safe/onboard-services?src=obs1"
$sacru.Range("A7").Value = "https://www.delta.com/us/en/travel-update-center/flying-what-you-need-to-know/coronavirus-regional-restrictions?src=restrictions1"
$sacru.Range("A8").Value = "https://www.delta.com/us/en/travel-update-center/from-delta-to-our-customers/medical-volunteers-book-free?src=flyfree1"
$sacru.Range("A9").Value = "https://www.delta.com/us/en/delta-digital/mobile?src=app1"
$sacru.Range("A10").Value = "https://www.delta.com/us/en/coronavirus-update-center/ways-we-are-keeping-you-safe/skyclub-update?src=skyclub1"
$sacru.Columns.Item(1).ColumnWidth = 124.30729166666667

# ---------------------------------------------------------------------------
# 2) OUATT - Header Names
# ---------------------------------------------------------------------------
$ouattHeader = Add-SheetAtEnd "OUATT - Header Names"
$ouattHeader.Range("A2").Value = "DOWNLOAD THE FLY DELTA APP"
$ouattHeader.Range("A3").Value = "FLYING WITH A PARTNER AIRLINE?"
$ouattHeader.Range("A1").Value = "DELTA SKY CLUB UPDATES"
$ouattHeader.Columns.Item(1).ColumnWidth = 30.736979166666668
$ouattHeader.Activate()
$ouattHeader.Range("C5").Select()
$ouattHeader.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3) OUATT - DeltaSkyClub URL
# ---------------------------------------------------------------------------
$ouattSkyClub = Add-SheetAtEnd "OUATT - DeltaSkyClub URL"
$ouattSkyClub.Range("A1").Value = "https://www.delta.com/us/en/travel-update-center/ways-we-are-keeping-you-safe/skyclub-update"
$ouattSkyClub.Columns.Item(1).ColumnWidth = 92.45182291666667

# ---------------------------------------------------------------------------
# 4) OUATT - FlyDeltaApp URL
# ---------------------------------------------------------------------------
$ouattFlyDelta = Add-SheetAtEnd "OUATT - FlyDeltaApp URL"
$ouattFlyDelta.Range("A1").Value = "https://www.delta.com/us/en/delta-digital/mobile?src=app2"
$ouattFlyDelta.Columns.Item(1).ColumnWidth = 56.307291666666664
$ouattFlyDelta.Activate()
$ouattFlyDelta.Range("F24").Select()

# ---------------------------------------------------------------------------
# 5) OUATT - FlyingPartnerAir URL
# ---------------------------------------------------------------------------
$ouattPartnerAir = Add-SheetAtEnd "OUATT - FlyingPartnerAir URL"
$ouattPartnerAir.Range("A1").Value = "https://www.delta.com/us/en/travel-update-center/flying-what-you-need-to-know/things-to-know-when-you-travel-with-a-partner-airline"

# ---------------------------------------------------------------------------
# 6) FAQ Header Names
# ---------------------------------------------------------------------------
$faqHeader = Add-SheetAtEnd "FAQ Header Names"
$faqHeader.Range("A1").Value = "TRAVEL FLEXIBILITY"
$faqHeader.Range("A2").Value = "SAFER TRAVEL"
$faqHeader.Range("A3").Value = "DELTA PARTNERS AND AFFILIATES"
$faqHeader.Columns.Item(1).ColumnWidth = 30.260625

# Make the new last sheet the active / selected tab, matching the target
# workbook view (activeTab moves from "CUC - TravelingWithUs Names" to the
# newly appended "FAQ Header Names" sheet).
$faqHeader.Activate()
$faqHeader.Range("G3").Select()
